$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_LT"

$values = @(
    85.022385022385023,
    85.714285714285708,
    86.080586080586087,
    88.970288970288962,
    89.662189662189661,
    89.417989417989418,
    80.138380138380143,
    81.277981277981283,
    81.03378103378104,
    80.545380545380539,
    79.446479446479444,
    81.847781847781846,
    80.219780219780219,
    80.504680504680508,
    82.4989824989825,
    79.93487993487993,
    81.522181522181526,
    86.853886853886848,
    90.476190476190482,
    90.923890923890923,
    90.516890516890513,
    83.597883597883595,
    86.65038665038665,
    86.080586080586087,
    83.801383801383793,
    83.068783068783063,
    83.597883597883595,
    81.400081400081405,
    80.708180708180706,
    80.341880341880341,
    89.865689865689873,
    93.121693121693113,
    92.673992673992672,
    87.179487179487182,
    91.656491656491653,
    82.539682539682531,
    87.993487993487989,
    85.144485144485145,
    84.167684167684172,
    80.504680504680508,
    81.725681725681724,
    81.359381359381359,
    81.359381359381359,
    84.126984126984127,
    85.592185592185587,
    79.812779812779809,
    81.196581196581192,
    81.03378103378104
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
